# The two slides currently use the "Blank" layout and carry no shapes
# of their own. The edit re-applies the "Title and Content" custom
# layout to each slide (this materializes the inherited Title and
# Content placeholder shapes on the slide itself) and then fills the
# placeholders in with the same outline prompt text / levels that the
# "Title and Content" layout defines.

$p = $ppt.ActivePresentation
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # Resetting the slide's layout materializes the Title (id 2) and
    # Content Placeholder (id 3) shapes on the slide.
    $s.CustomLayout = $titleAndContent

    $title = $s.Shapes.Item(1)
    $title.TextFrame.TextRange.Text = "Click to edit Master title style"

    $body = $s.Shapes.Item(2)
    $bodyText = $body.TextFrame.TextRange
    $bodyText.Text = "Click to edit Master text styles"
    $null = $bodyText.InsertAfter("`rSecond level")
    $null = $bodyText.InsertAfter("`rThird level")
    $null = $bodyText.InsertAfter("`rFourth level")
    $null = $bodyText.InsertAfter("`rFifth level")

    $bodyText.Lines(1, 1).IndentLevel = 1
    $bodyText.Lines(2, 1).IndentLevel = 2
    $bodyText.Lines(3, 1).IndentLevel = 3
    $bodyText.Lines(4, 1).IndentLevel = 4
    $bodyText.Lines(5, 1).IndentLevel = 5
}
